$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Add the two new timeline rows (5 and 6), mirroring the formatting of
#    the existing data row (row 4).
# ---------------------------------------------------------------------------

# Pre-fill values first (so autofit / paste-format operations don't need to
# recompute anything based on still-empty cells).
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 45508
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "Made login/signup page with some login functionality"

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "30/3/2024"
$ws.Range("C6").Value = 3.25
$ws.Range("D6").Value = "Finished Login/SignUp/Forgot Password functionality in UI"

# Copy the formatting of row 4 down onto the two new rows.
$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A5:D6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Move the "Total hours Spent" block from rows 7-9 down to rows 10-12.
# ---------------------------------------------------------------------------

# Set the new totals content first (avoids triggering an autofit row-height
# recalculation once the bold/large font style is applied).
$ws.Range("C10").Value = "Total hours Spent"
$ws.Range("D10").Formula = "=SUM(C4:C9)"

# Break up the old merges so the source cells can be copied individually.
$ws.Range("C7:C9").UnMerge() | Out-Null
$ws.Range("D7:D9").UnMerge() | Out-Null

# Carry the formatting from the old block onto the new block.
$ws.Range("C7:D9").Copy() | Out-Null
$ws.Range("C10:D12").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

# Clear the old block's contents (formula/value/text) now that it has moved.
$ws.Range("C7:D9").ClearContents() | Out-Null

# Re-merge at the new location.
$ws.Range("C10:C12").Merge() | Out-Null
$ws.Range("D10:D12").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 3) Fill in the (blank) formatting-only cells left behind on rows 7-9.
# ---------------------------------------------------------------------------

$ws.Range("A7:D8").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null      # xlPasteFormats (keeps A7 as-is / General)
$excel.CutCopyMode = 0

$ws.Range("A4").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B7").NumberFormat = "mm-dd-yy"
$ws.Range("B9").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# 4) Misc view state to mirror the author's final selection.
# ---------------------------------------------------------------------------
$ws.Range("D20").Select() | Out-Null

$wb.Save()
